# Tareas a realizar Practica evaluacion.docx
#
# The only textual/content change in this revision is that the name
# "Adrian" was filled in as the 4th team member on the
# "4ª persona():" line, turning it into "4ª persona(Adrian):".
# (Word's proofing engine also re-ran in the background on this pass,
# which is why the canonical XML shows extra <w:proofErr/> markers and
# some runs split at spell/grammar-check boundaries -- none of that
# changes the visible text.)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "4ª persona():",  # FindText
    $true,             # MatchCase
    $false,            # MatchWholeWord
    $false,            # MatchWildcards
    $false,            # MatchSoundsLike
    $false,            # MatchAllWordForms
    $true,             # Forward
    1,                 # Wrap (wdFindContinue)
    $false,            # Format
    "4ª persona(Adrian):",  # ReplaceWith
    2                  # Replace (wdReplaceAll)
)
